$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "01-introduction.html"
$ws.Range("E3").Value = "02-mles"
$ws.Range("C3").Value = "Lecture 2: MLEs & Projections"
$ws.Range("D3").Value = "02-MLEs.html"
$ws.Range("F4").Value = "lab01.html"
